$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume figures (plus a few coin
# name/link re-rankings) as captured by the latest GitHub Actions run.
# Numeric-looking Price values get a leading apostrophe so Excel keeps
# storing them as literal text (e.g. "598.66"), matching the original
# inline-string cells instead of converting them into real numbers.

$ws.Range("D2").Value = "65.885.37"
$ws.Range("D3").Value = "2.664.89"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'598.66"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").Value = "'158.74"
$ws.Range("E6").Value = "  +1.16%  "
$ws.Range("D7").Value = "'0.651"
$ws.Range("E7").Value = "  +4.72%  "
$ws.Range("E9").Value = "  -2.18%  "
$ws.Range("D10").Value = "'0.402"
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("E12").Value = "  +1.65%  "
$ws.Range("D13").Value = "'29.09"
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("D14").Value = "'0.0000196"
$ws.Range("E14").Value = "  -1.65%  "
$ws.Range("D15").Value = "3.142.85"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("D16").Value = "65.732.38"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "2.655.94"
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").Value = "'12.64"
$ws.Range("E18").Value = "  -2.06%  "
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'352.48"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'7.51"
$ws.Range("E21").Value = "  -1.37%  "
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").Value = "'69.94"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").Value = "'1.83"
$ws.Range("E24").Value = "  +11.31%  "
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("D26").Value = "'9.67"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").Value = "  +1.59%  "
$ws.Range("D28").Value = "'575.78"
$ws.Range("E28").Value = "  +8.75%  "
$ws.Range("D29").Value = "'8.21"
$ws.Range("E29").Value = "  +1.59%  "
$ws.Range("E30").Value = "  -2.54%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'2.16"
$ws.Range("E31").Value = "  +1.00%  "
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("E33").Value = "  +3.12%  "
$ws.Range("D34").Value = "'6.77"
$ws.Range("E34").Value = "  +4.45%  "
$ws.Range("E35").Value = "  +1.14%  "
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").Value = "'20.63"
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("E39").Value = "  +1.23%  "
$ws.Range("D40").Value = "'154.53"
$ws.Range("E40").Value = "  -2.19%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'161.85"
$ws.Range("E41").Value = "  -1.50%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'4.11"
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.33"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("D44").Value = "'0.0618"
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'23.27"
$ws.Range("E45").Value = "  +1.46%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.646"
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0258"
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.103"
$ws.Range("E48").Value = "  +1.90%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'19.85"
$ws.Range("E49").Value = "  -1.65%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0247"
$ws.Range("E50").Value = "  -7.92%  "
$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").Value = "'0.815"
$ws.Range("E51").Value = "  -0.31%  "
